$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing sd_total value for the existing Apr 12 row (39)
$ws.Range("B39").Value = 898

# New row 40 - Apr 13
$ws.Range("A40").Value = 43934
$ws.Range("B40").Value = 919
$ws.Range("C40").Value = 1847
$ws.Range("D40").Value = 11
$ws.Range("F40").Value = 23
$ws.Range("H40").Value = 279
$ws.Range("J40").Value = 359
$ws.Range("L40").Value = 320
$ws.Range("N40").Value = 338
$ws.Range("P40").Value = 255
$ws.Range("R40").Value = 146
$ws.Range("T40").Value = 113
$ws.Range("V40").Value = 3
$ws.Range("X40").Value = 905
$ws.Range("Y40").Value = 937
$ws.Range("Z40").Value = 5
$ws.Range("AA40").Value = 420
$ws.Range("AB40").Value = 156
$ws.Range("AC40").Value = 47

# New row 41 - Apr 14
$ws.Range("A41").Value = 43935
$ws.Range("C41").Value = 1930
$ws.Range("D41").Value = 12
$ws.Range("F41").Value = 24
$ws.Range("H41").Value = 292
$ws.Range("J41").Value = 367
$ws.Range("L41").Value = 339
$ws.Range("N41").Value = 352
$ws.Range("P41").Value = 265
$ws.Range("R41").Value = 155
$ws.Range("T41").Value = 121
$ws.Range("V41").Value = 3
$ws.Range("X41").Value = 954
$ws.Range("Y41").Value = 971
$ws.Range("Z41").Value = 5
$ws.Range("AA41").Value = 450
$ws.Range("AB41").Value = 164
$ws.Range("AC41").Value = 53

# Update the selection to match the new active cell after data entry
$ws.Range("B41").Select() | Out-Null
